$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = "Bitcoin"
$data[0,1] = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$data[0,2] = "'28.274.81"
$data[0,3] = "'  +1.18%  "
$data[1,0] = "Ethereum"
$data[1,1] = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$data[1,2] = "'1.884.99"
$data[1,3] = "'  +1.55%  "
$data[2,0] = "TetherUSD"
$data[2,1] = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$data[2,2] = "'1.008"
$data[2,3] = "'  +0.32%  "
$data[3,0] = "BNB"
$data[3,1] = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$data[3,2] = "'314.46"
$data[3,3] = "'  +0.96%  "
$data[4,0] = "USDC"
$data[4,1] = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$data[4,2] = "'1.007"
$data[4,3] = "'  +0.36%  "
$data[5,0] = "XRP"
$data[5,1] = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$data[5,2] = "'0.5144"
$data[5,3] = "'  +0.55%  "
$data[6,0] = "Cardano"
$data[6,1] = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$data[6,2] = "'0.3920"
$data[6,3] = "'  +3.39%  "
$data[7,0] = "Dogecoin"
$data[7,1] = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$data[7,2] = "'0.08370"
$data[7,3] = "'  +0.25%  "
$data[8,0] = "Polygon"
$data[8,1] = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$data[8,2] = "'1.123"
$data[8,3] = "'  +1.76%  "
$data[9,0] = "OKB"
$data[9,1] = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$data[9,2] = "'41.64"
$data[9,3] = "'  +0.87%  "
$data[10,0] = "Polkadot"
$data[10,1] = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$data[10,2] = "'6.244"
$data[10,3] = "'  +0.95%  "
$data[11,0] = "WrappedEther"
$data[11,1] = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$data[11,2] = "'1.894.14"
$data[11,3] = "'  +1.80%  "
$data[12,0] = "Solana"
$data[12,1] = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$data[12,2] = "'20.75"
$data[12,3] = "'  +1.83%  "
$data[13,0] = "Chainlink"
$data[13,1] = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$data[13,2] = "'7.298"
$data[13,3] = "'  +2.07%  "
$data[14,0] = "BinanceUSD"
$data[14,1] = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$data[14,2] = "'1.007"
$data[14,3] = "'  +0.23%  "
$data[15,0] = "ShibaInu"
$data[15,1] = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$data[15,2] = "'0.00001110"
$data[15,3] = "'  +1.84%  "
$data[16,0] = "Litecoin"
$data[16,1] = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$data[16,2] = "'91.62"
$data[16,3] = "'  +1.78%  "
$data[17,0] = "TRON"
$data[17,1] = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$data[17,2] = "'0.06678"
$data[17,3] = "'  +0.80%  "
$data[18,0] = "Avalanche"
$data[18,1] = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$data[18,2] = "'17.82"
$data[18,3] = "'  +0.37%  "
$data[19,0] = "Dai"
$data[19,1] = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$data[19,2] = "'1.007"
$data[19,3] = "'  +0.40%  "
$data[20,0] = "Uniswap"
$data[20,1] = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$data[20,2] = "'6.067"
$data[20,3] = "'  +1.22%  "
$data[21,0] = "WrappedBTC"
$data[21,1] = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$data[21,2] = "'28.309.95"
$data[21,3] = "'  +1.16%  "
$data[22,0] = "Cosmos"
$data[22,1] = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$data[22,2] = "'11.17"
$data[22,3] = "'  +1.25%  "
$data[23,0] = "Toncoin"
$data[23,1] = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$data[23,2] = "'2.279"
$data[23,3] = "'  +1.14%  "
$data[24,0] = "WrappedliquidstakedEther2.0"
$data[24,1] = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$data[24,2] = "'2.097.39"
$data[24,3] = "'  +1.00%  "
$data[25,0] = "LidoDAOToken"
$data[25,1] = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$data[25,2] = "'2.525"
$data[25,3] = "'  -1.00%  "
$data[26,0] = "Monero"
$data[26,1] = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$data[26,2] = "'159.21"
$data[26,3] = "'  +1.38%  "
$data[27,0] = "EthereumClassic"
$data[27,1] = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$data[27,2] = "'20.69"
$data[27,3] = "'  +1.42%  "
$data[28,0] = "BitcoinCash"
$data[28,1] = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$data[28,2] = "'125.74"
$data[28,3] = "'  +0.27%  "
$data[29,0] = "Stellar"
$data[29,1] = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$data[29,2] = "'0.1067"
$data[29,3] = "'  +1.28%  "
$data[30,0] = "ImmutableX"
$data[30,1] = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$data[30,2] = "'1.051"
$data[30,3] = "'  +1.50%  "
$data[31,0] = "Filecoin"
$data[31,1] = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$data[31,2] = "'5.921"
$data[31,3] = "'  +6.65%  "
$data[32,0] = "HuobiToken"
$data[32,1] = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$data[32,2] = "'3.609"
$data[32,3] = "'  +0.20%  "
$data[33,0] = "FraxShare"
$data[33,1] = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$data[33,2] = "'9.823"
$data[33,3] = "'  +3.25%  "
$data[34,0] = "VeChain"
$data[34,1] = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$data[34,2] = "'0.02469"
$data[34,3] = "'  +2.88%  "
$data[35,0] = "Hedera"
$data[35,1] = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$data[35,2] = "'0.06590"
$data[35,3] = "'  +1.56%  "
$data[36,0] = "Algorand"
$data[36,1] = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$data[36,2] = "'0.2199"
$data[36,3] = "'  +2.58%  "
$data[37,0] = "ARBITRUM"
$data[37,1] = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$data[37,2] = "'1.213"
$data[37,3] = "'  +0.78%  "
$data[38,0] = "TheSandbox"
$data[38,1] = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$data[38,2] = "'0.6562"
$data[38,3] = "'  +3.41%  "
$data[39,0] = "InternetComputer(DFINITY)"
$data[39,1] = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$data[39,2] = "'5.033"
$data[39,3] = "'  +4.03%  "
$data[40,0] = "TrustWalletToken"
$data[40,1] = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$data[40,2] = "'1.233"
$data[40,3] = "'  +0.67%  "
$data[41,0] = "Aptos"
$data[41,1] = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$data[41,2] = "'11.32"
$data[41,3] = "'  +0.87%  "
$data[42,0] = "Decentraland"
$data[42,1] = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$data[42,2] = "'0.6160"
$data[42,3] = "'  +2.82%  "
$data[43,0] = "EnergySwap"
$data[43,1] = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$data[43,2] = "'13.13"
$data[43,3] = "'  +1.88%  "
$data[44,0] = "WEMIXTOKEN"
$data[44,1] = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$data[44,2] = "'1.289"
$data[44,3] = "'  +0.47%  "
$data[45,0] = "PancakeSwap"
$data[45,1] = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$data[45,2] = "'3.683"
$data[45,3] = "'  +0.72%  "
$data[46,0] = "NEARProtocol"
$data[46,1] = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$data[46,2] = "'2.019"
$data[46,3] = "'  +2.69%  "
$data[47,0] = "EOS"
$data[47,1] = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$data[47,2] = "'1.238"
$data[47,3] = "'  +2.86%  "
$data[48,0] = "Quant"
$data[48,1] = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$data[48,2] = "'121.66"
$data[48,3] = "'  +1.02%  "
$data[49,0] = "Aave"
$data[49,1] = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$data[49,2] = "'79.18"
$data[49,3] = "'  -0.20%  "

$ws.Range("B2:E51").Value = $data